$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.861.88"

$ws.Range("D3").Value = "1.880.43"
$ws.Range("E3").Value = "  -0.36%  "

$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'336.25"
$ws.Range("E5").Value = "  +0.59%  "

$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("D7").Value = "'0.4704"
$ws.Range("E7").Value = "  -0.45%  "

$ws.Range("D8").Value = "'0.3951"
$ws.Range("E8").Value = "  +0.41%  "

$ws.Range("D9").Value = "'45.69"
$ws.Range("E9").Value = "  -4.11%  "

$ws.Range("D10").Value = "'0.08034"
$ws.Range("E10").Value = "  -0.68%  "

$ws.Range("D11").Value = "'1.015"
$ws.Range("E11").Value = "  -1.29%  "

$ws.Range("D12").Value = "'22.09"
$ws.Range("E12").Value = "  -0.42%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'6.019"
$ws.Range("E13").Value = "  +0.48%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.865.15"
$ws.Range("E14").Value = "  -1.11%  "

$ws.Range("D15").Value = "'7.303"
$ws.Range("E15").Value = "  +2.27%  "

$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("D17").Value = "'89.18"
$ws.Range("E17").Value = "  +2.08%  "

$ws.Range("D18").Value = "'0.06728"
$ws.Range("E18").Value = "  -0.12%  "

$ws.Range("D19").Value = "'0.00001046"
$ws.Range("E19").Value = "  -0.52%  "

$ws.Range("D20").Value = "'17.36"
$ws.Range("E20").Value = "  -0.10%  "

$ws.Range("D21").Value = "'1.011"
$ws.Range("E21").Value = "  +0.42%  "

$ws.Range("D22").Value = "27.851.12"
$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D23").Value = "'5.507"
$ws.Range("E23").Value = "  -0.42%  "

$ws.Range("D24").Value = "'11.03"
$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("D25").Value = "'2.317"
$ws.Range("E25").Value = "  -0.55%  "

$ws.Range("D26").Value = "2.093.85"
$ws.Range("E26").Value = "  -0.79%  "

$ws.Range("D27").Value = "'159.52"
$ws.Range("E27").Value = "  +0.26%  "

$ws.Range("D28").Value = "'19.89"
$ws.Range("E28").Value = "  -1.57%  "

$ws.Range("D29").Value = "'2.166"
$ws.Range("E29").Value = "  +2.72%  "

$ws.Range("D30").Value = "'5.504"
$ws.Range("E30").Value = "  -1.41%  "

$ws.Range("D31").Value = "'122.16"
$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("D32").Value = "'0.9893"
$ws.Range("E32").Value = "  +0.58%  "

$ws.Range("D33").Value = "'0.09538"
$ws.Range("E33").Value = "  +0.39%  "

$ws.Range("D34").Value = "'3.634"
$ws.Range("E34").Value = "  +0.41%  "

$ws.Range("D35").Value = "'5.357"
$ws.Range("E35").Value = "  -0.08%  "

$ws.Range("D36").Value = "'1.355"
$ws.Range("E36").Value = "  -6.74%  "

$ws.Range("D37").Value = "'0.06085"
$ws.Range("E37").Value = "  -1.28%  "

$ws.Range("D38").Value = "'0.02249"
$ws.Range("E38").Value = "  -0.79%  "

$ws.Range("D39").Value = "'1.203"
$ws.Range("E39").Value = "  -1.58%  "

$ws.Range("D40").Value = "'8.341"
$ws.Range("E40").Value = "  +3.01%  "

$ws.Range("D41").Value = "'1.007"
$ws.Range("E41").Value = "  +0.08%  "

$ws.Range("D42").Value = "'0.6015"
$ws.Range("E42").Value = "  +0.04%  "

$ws.Range("D43").Value = "'0.1898"
$ws.Range("E43").Value = "  -0.04%  "

$ws.Range("E44").Value = "  +1.16%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.5693"
$ws.Range("E45").Value = "  -0.48%  "

$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.248"
$ws.Range("E46").Value = "  -0.95%  "

$ws.Range("D47").Value = "'12.24"
$ws.Range("E47").Value = "  +0.03%  "

$ws.Range("D48").Value = "'1.946"
$ws.Range("E48").Value = "  -0.15%  "

$ws.Range("E49").Value = "  -2.00%  "

$ws.Range("D50").Value = "'112.63"
$ws.Range("E50").Value = "  -1.36%  "

$ws.Range("E51").Value = "  -10.44%  "

# Strip the quote-prefix formatting introduced by the leading apostrophe trick above,
# restoring cells to their original (unstyled) appearance while keeping text-typed values.
$clearRanges = @("D4","D5","D7","D8","D9","D10","D11","D12","D13","D15","D17","D18","D19","D20","D21","D23","D24","D25","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D45","D46","D47","D48","D50")
foreach ($addr in $clearRanges) {
    $ws.Range($addr).ClearFormats()
}
